# Append: 2025-11-30 06:35 JST
# Update the "取得日時" (retrieved timestamp) column on the "ランサーズ" sheet
# for all existing data rows (2-9) from the old scrape timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-30 06:35:12"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
